$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.336.88'
$ws.Range("E2").Value = '  -2.32%  '

$ws.Range("D3").Value = '1.794.02'
$ws.Range("E3").Value = '  -2.05%  '

$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.32%  '

$ws.Range("D6").Value = '306.85'
$ws.Range("E6").Value = '  -1.40%  '

$ws.Range("D7").Value = '0.4503'
$ws.Range("E7").Value = '  -1.63%  '

$ws.Range("D8").Value = '0.3594'
$ws.Range("E8").Value = '  -3.07%  '

$ws.Range("D9").Value = '45.84'
$ws.Range("E9").Value = '  -0.30%  '

$ws.Range("E10").Value = '  -1.40%  '

$ws.Range("D11").Value = '0.8840'
$ws.Range("E11").Value = '  +0.67%  '

$ws.Range("E12").Value = '  -1.01%  '

$ws.Range("D13").Value = '19.39'
$ws.Range("E13").Value = '  -1.47%  '

$ws.Range("D14").Value = '1.823.06'
$ws.Range("E14").Value = '  -0.57%  '

$ws.Range("D15").Value = '5.278'
$ws.Range("E15").Value = '  -1.21%  '

$ws.Range("D16").Value = '6.325'
$ws.Range("E16").Value = '  -1.36%  '

$ws.Range("D17").Value = '84.82'
$ws.Range("E17").Value = '  -2.86%  '

$ws.Range("E18").Value = '  -0.16%  '

$ws.Range("D19").Value = '0.000008504'
$ws.Range("E19").Value = '  -2.52%  '

$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  -0.20%  '

$ws.Range("E21").Value = '  -1.82%  '

$ws.Range("D22").Value = '26.367.88'
$ws.Range("E22").Value = '  -2.28%  '

$ws.Range("D23").Value = '4.972'
$ws.Range("E23").Value = '  -0.82%  '

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '10.52'
$ws.Range("E24").Value = '  +0.76%  '

$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value = '2.013.81'
$ws.Range("E25").Value = '  -2.11%  '

$ws.Range("D26").Value = '1.973'
$ws.Range("E26").Value = '  -2.73%  '

$ws.Range("D27").Value = '151.40'
$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("D28").Value = '17.82'
$ws.Range("E28").Value = '  -2.31%  '

$ws.Range("D29").Value = '2.020'
$ws.Range("E29").Value = '  +2.37%  '

$ws.Range("D30").Value = '111.77'
$ws.Range("E30").Value = '  -2.05%  '

$ws.Range("D31").Value = '4.886'
$ws.Range("E31").Value = '  -1.23%  '

$ws.Range("D32").Value = '0.08671'
$ws.Range("E32").Value = '  -1.55%  '

$ws.Range("D33").Value = '3.061'
$ws.Range("E33").Value = '  +1.00%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '4.442'
$ws.Range("E34").Value = '  -0.93%  '

$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").Value = '2.726'
$ws.Range("E35").Value = '  +5.88%  '

$ws.Range("D36").Value = '0.7236'
$ws.Range("E36").Value = '  -4.16%  '

$ws.Range("D37").Value = '1.105'
$ws.Range("E37").Value = '  -2.87%  '

$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("D39").Value = '1.066'
$ws.Range("E39").Value = '  -2.24%  '

$ws.Range("D40").Value = '0.01929'
$ws.Range("E40").Value = '  -0.57%  '

$ws.Range("D41").Value = '0.05090'
$ws.Range("E41").Value = '  -1.26%  '

$ws.Range("E42").Value = '  -1.11%  '

$ws.Range("D43").Value = '0.5074'
$ws.Range("E43").Value = '  +1.70%  '

$ws.Range("D44").Value = '6.858'
$ws.Range("E44").Value = '  -1.55%  '

$ws.Range("D45").Value = '0.1515'
$ws.Range("E45").Value = '  -5.37%  '

$ws.Range("D46").Value = '7.996'
$ws.Range("E46").Value = '  -4.13%  '

$ws.Range("E47").Value = '  -0.24%  '

$ws.Range("D48").Value = '0.4632'
$ws.Range("E48").Value = '  -1.34%  '

$ws.Range("D49").Value = '101.08'
$ws.Range("E49").Value = '  -1.31%  '

$ws.Range("D50").Value = '9.825'
$ws.Range("E50").Value = '  -3.06%  '

$ws.Range("D51").Value = '1.578'
$ws.Range("E51").Value = '  -2.36%  '
